$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = 1414.76449725
$wsSchedule.Range("F2").Value = 23.3922701264881

# --- Sheet "Detailed" ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B17").Value = 51.40072
$wsDetailed.Range("B18").Value = 50.14948

$wsDetailed.Range("B19").Value = 36.06
$wsDetailed.Range("C19").Value = "historical"

$wsDetailed.Range("B20").Value = 36.06
$wsDetailed.Range("C20").Value = "historical"

$wsDetailed.Range("B21").Value = 0.08382000000000001
$wsDetailed.Range("B22").Value = 0.02167
$wsDetailed.Range("B23").Value = 22.07
$wsDetailed.Range("B24").Value = -0.0001
$wsDetailed.Range("B25").Value = -2.52488
$wsDetailed.Range("B26").Value = 0.51
$wsDetailed.Range("B27").Value = 0.51
$wsDetailed.Range("B28").Value = -4.66156
$wsDetailed.Range("B29").Value = -5.50985
$wsDetailed.Range("B30").Value = 0.51
$wsDetailed.Range("B31").Value = 36.0601
$wsDetailed.Range("B32").Value = 36.0601

$wsDetailed.Range("B34").Value = 40.61245
$wsDetailed.Range("B35").Value = 44.01368
$wsDetailed.Range("B36").Value = 47.46912
$wsDetailed.Range("B37").Value = 19.05091
$wsDetailed.Range("B38").Value = 49.58561
$wsDetailed.Range("B39").Value = 70.46706

$wsDetailed.Range("B41").Value = 120.01
$wsDetailed.Range("B42").Value = 140.37518

$wsDetailed.Range("B44").Value = 105.79004
$wsDetailed.Range("B45").Value = 105.79
$wsDetailed.Range("B46").Value = 85.95
$wsDetailed.Range("B47").Value = 65.31265
$wsDetailed.Range("B48").Value = 67.74731
$wsDetailed.Range("B49").Value = 74.51801
